# Update NATMI ligand/receptor TPM-derived statistics for Gas6-Axl.
# Rewrites the numeric columns (G..T, excluding the unchanged K/L counts)
# on rows 2-26 of the active sheet with the values recomputed from the new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 16.014007
$ws.Range("H2").Value = 48.04202100000001
$ws.Range("I2").Value = 0.09359269140871758
$ws.Range("J2").Value = 0.09359269140871758
$ws.Range("M2").Value = 3.243447333333334
$ws.Range("N2").Value = 9.730342
$ws.Range("O2").Value = 0.01255381554477167
$ws.Range("P2").Value = 0.01255381554477167
$ws.Range("Q2").Value = 51.94058830013135
$ws.Range("R2").Value = 467.4652947011821
$ws.Range("S2").Value = 0.001174945384283777
$ws.Range("T2").Value = 0.001174945384283776

# Row 3
$ws.Range("G3").Value = 16.014007
$ws.Range("H3").Value = 48.04202100000001
$ws.Range("I3").Value = 0.09359269140871758
$ws.Range("J3").Value = 0.09359269140871758
$ws.Range("O3").Value = 0.2470858318854151
$ws.Range("P3").Value = 0.2470858318854151
$ws.Range("Q3").Value = 1022.301420869669
$ws.Range("R3").Value = 9200.712787827018
$ws.Range("S3").Value = 0.02312542801511793
$ws.Range("T3").Value = 0.02312542801511793

# Row 4
$ws.Range("G4").Value = 16.014007
$ws.Range("H4").Value = 48.04202100000001
$ws.Range("I4").Value = 0.09359269140871758
$ws.Range("J4").Value = 0.09359269140871758
$ws.Range("M4").Value = 103.7552793333333
$ws.Range("N4").Value = 311.265838
$ws.Range("O4").Value = 0.4015864926064037
$ws.Range("P4").Value = 0.4015864926064037
$ws.Range("Q4").Value = 1661.537769530956
$ws.Range("R4").Value = 14953.8399257786
$ws.Range("S4").Value = 0.03758556067642039
$ws.Range("T4").Value = 0.03758556067642038

# Row 5
$ws.Range("G5").Value = 16.014007
$ws.Range("H5").Value = 48.04202100000001
$ws.Range("I5").Value = 0.09359269140871758
$ws.Range("J5").Value = 0.09359269140871758
$ws.Range("M5").Value = 23.61337433333334
$ws.Range("N5").Value = 70.84012300000001
$ws.Range("O5").Value = 0.09139594860190291
$ws.Range("P5").Value = 0.09139594860190289
$ws.Range("Q5").Value = 378.1447418676204
$ws.Range("R5").Value = 3403.302676808584
$ws.Range("S5").Value = 0.008553992813504911
$ws.Range("T5").Value = 0.008553992813504911

# Row 6
$ws.Range("G6").Value = 16.014007
$ws.Range("H6").Value = 48.04202100000001
$ws.Range("I6").Value = 0.09359269140871758
$ws.Range("J6").Value = 0.09359269140871758
$ws.Range("M6").Value = 63.91341533333334
$ws.Range("N6").Value = 191.740246
$ws.Range("O6").Value = 0.2473779113615065
$ws.Range("P6").Value = 0.2473779113615065
$ws.Range("Q6").Value = 1023.509880541908
$ws.Range("R6").Value = 9211.588924877167
$ws.Range("S6").Value = 0.02315276451939057
$ws.Range("T6").Value = 0.02315276451939056

# Row 7
$ws.Range("I7").Value = 0.5957388935007043
$ws.Range("J7").Value = 0.5957388935007044
$ws.Range("M7").Value = 3.243447333333334
$ws.Range("N7").Value = 9.730342
$ws.Range("O7").Value = 0.01255381554477167
$ws.Range("P7").Value = 0.01255381554477167
$ws.Range("Q7").Value = 330.6137278024011
$ws.Range("R7").Value = 2975.52355022161
$ws.Range("S7").Value = 0.007478796181854216
$ws.Range("T7").Value = 0.007478796181854217

# Row 8
$ws.Range("I8").Value = 0.5957388935007043
$ws.Range("J8").Value = 0.5957388935007044
$ws.Range("O8").Value = 0.2470858318854151
$ws.Range("P8").Value = 0.2470858318854151
$ws.Range("Q8").Value = 6507.182431943265
$ws.Range("S8").Value = 0.1471986400871182
$ws.Range("T8").Value = 0.1471986400871183

# Row 9
$ws.Range("I9").Value = 0.5957388935007043
$ws.Range("J9").Value = 0.5957388935007044
$ws.Range("M9").Value = 103.7552793333333
$ws.Range("N9").Value = 311.265838
$ws.Range("O9").Value = 0.4015864926064037
$ws.Range("P9").Value = 0.4015864926064037
$ws.Range("Q9").Value = 10576.06803940892
$ws.Range("R9").Value = 95184.61235468028
$ws.Range("S9").Value = 0.2392406927501677
$ws.Range("T9").Value = 0.2392406927501677

# Row 10
$ws.Range("I10").Value = 0.5957388935007043
$ws.Range("J10").Value = 0.5957388935007044
$ws.Range("M10").Value = 23.61337433333334
$ws.Range("N10").Value = 70.84012300000001
$ws.Range("O10").Value = 0.09139594860190291
$ws.Range("P10").Value = 0.09139594860190289
$ws.Range("Q10").Value = 2406.977796156663
$ws.Range("R10").Value = 21662.80016540996
$ws.Range("S10").Value = 0.05444812129054488
$ws.Range("T10").Value = 0.05444812129054488

# Row 11
$ws.Range("I11").Value = 0.5957388935007043
$ws.Range("J11").Value = 0.5957388935007044
$ws.Range("M11").Value = 63.91341533333334
$ws.Range("N11").Value = 191.740246
$ws.Range("O11").Value = 0.2473779113615065
$ws.Range("P11").Value = 0.2473779113615065
$ws.Range("Q11").Value = 6514.874554235548
$ws.Range("R11").Value = 58633.87098811993
$ws.Range("S11").Value = 0.1473726431910192
$ws.Range("T11").Value = 0.1473726431910192

# Row 12
$ws.Range("G12").Value = 19.33193133333333
$ws.Range("H12").Value = 57.995794
$ws.Range("I12").Value = 0.1129840572453343
$ws.Range("J12").Value = 0.1129840572453343
$ws.Range("M12").Value = 3.243447333333334
$ws.Range("N12").Value = 9.730342
$ws.Range("O12").Value = 0.01255381554477167
$ws.Range("P12").Value = 0.01255381554477167
$ws.Range("Q12").Value = 62.70210113128311
$ws.Range("R12").Value = 564.3189101815479
$ws.Range("S12").Value = 0.00141838101415785
$ws.Range("T12").Value = 0.00141838101415785

# Row 13
$ws.Range("G13").Value = 19.33193133333333
$ws.Range("H13").Value = 57.995794
$ws.Range("I13").Value = 0.1129840572453343
$ws.Range("J13").Value = 0.1129840572453343
$ws.Range("O13").Value = 0.2470858318854151
$ws.Range("P13").Value = 0.2470858318854151
$ws.Range("Q13").Value = 1234.110917412584
$ws.Range("R13").Value = 11106.99825671325
$ws.Range("S13").Value = 0.0279167597742528
$ws.Range("T13").Value = 0.02791675977425279

# Row 14
$ws.Range("G14").Value = 19.33193133333333
$ws.Range("H14").Value = 57.995794
$ws.Range("I14").Value = 0.1129840572453343
$ws.Range("J14").Value = 0.1129840572453343
$ws.Range("M14").Value = 103.7552793333333
$ws.Range("N14").Value = 311.265838
$ws.Range("O14").Value = 0.4015864926064037
$ws.Range("P14").Value = 0.4015864926064037
$ws.Range("Q14").Value = 2005.789935542819
$ws.Range("R14").Value = 18052.10941988537
$ws.Range("S14").Value = 0.04537287126959495
$ws.Range("T14").Value = 0.04537287126959494

# Row 15
$ws.Range("G15").Value = 19.33193133333333
$ws.Range("H15").Value = 57.995794
$ws.Range("I15").Value = 0.1129840572453343
$ws.Range("J15").Value = 0.1129840572453343
$ws.Range("M15").Value = 23.61337433333334
$ws.Range("N15").Value = 70.84012300000001
$ws.Range("O15").Value = 0.09139594860190291
$ws.Range("P15").Value = 0.09139594860190289
$ws.Range("Q15").Value = 456.4921311602959
$ws.Range("R15").Value = 4108.429180442662
$ws.Range("S15").Value = 0.01032628508882903
$ws.Range("T15").Value = 0.01032628508882903

# Row 16
$ws.Range("G16").Value = 19.33193133333333
$ws.Range("H16").Value = 57.995794
$ws.Range("I16").Value = 0.1129840572453343
$ws.Range("J16").Value = 0.1129840572453343
$ws.Range("M16").Value = 63.91341533333334
$ws.Range("N16").Value = 191.740246
$ws.Range("O16").Value = 0.2473779113615065
$ws.Range("P16").Value = 0.2473779113615065
$ws.Range("Q16").Value = 1235.569756502814
$ws.Range("R16").Value = 11120.12780852532
$ws.Range("S16").Value = 0.02794976009849969
$ws.Range("T16").Value = 0.02794976009849969

# Row 17
$ws.Range("G17").Value = 0.6875779999999999
$ws.Range("H17").Value = 2.062734
$ws.Range("I17").Value = 0.004018499278376935
$ws.Range("J17").Value = 0.004018499278376936
$ws.Range("M17").Value = 3.243447333333334
$ws.Range("N17").Value = 9.730342
$ws.Range("O17").Value = 0.01255381554477167
$ws.Range("P17").Value = 0.01255381554477167
$ws.Range("Q17").Value = 2.230123030558667
$ws.Range("R17").Value = 20.071107275028
$ws.Range("S17").Value = 0.00005044749870754211
$ws.Range("T17").Value = 0.00005044749870754211

# Row 18
$ws.Range("G18").Value = 0.6875779999999999
$ws.Range("H18").Value = 2.062734
$ws.Range("I18").Value = 0.004018499278376935
$ws.Range("J18").Value = 0.004018499278376936
$ws.Range("O18").Value = 0.2470858318854151
$ws.Range("P18").Value = 0.2470858318854151
$ws.Range("Q18").Value = 43.89357181864133
$ws.Range("R18").Value = 395.042146367772
$ws.Range("S18").Value = 0.0009929142371287055
$ws.Range("T18").Value = 0.0009929142371287055

# Row 19
$ws.Range("G19").Value = 0.6875779999999999
$ws.Range("H19").Value = 2.062734
$ws.Range("I19").Value = 0.004018499278376935
$ws.Range("J19").Value = 0.004018499278376936
$ws.Range("M19").Value = 103.7552793333333
$ws.Range("N19").Value = 311.265838
$ws.Range("O19").Value = 0.4015864926064037
$ws.Range("P19").Value = 0.4015864926064037
$ws.Range("Q19").Value = 71.33984745345465
$ws.Range("R19").Value = 642.0586270810919
$ws.Range("S19").Value = 0.001613775030744758
$ws.Range("T19").Value = 0.001613775030744758

# Row 20
$ws.Range("G20").Value = 0.6875779999999999
$ws.Range("H20").Value = 2.062734
$ws.Range("I20").Value = 0.004018499278376935
$ws.Range("J20").Value = 0.004018499278376936
$ws.Range("M20").Value = 23.61337433333334
$ws.Range("N20").Value = 70.84012300000001
$ws.Range("O20").Value = 0.09139594860190291
$ws.Range("P20").Value = 0.09139594860190289
$ws.Range("Q20").Value = 16.23603669736467
$ws.Range("R20").Value = 146.124330276282
$ws.Range("S20").Value = 0.0003672745535033223
$ws.Range("T20").Value = 0.0003672745535033223

# Row 21
$ws.Range("G21").Value = 0.6875779999999999
$ws.Range("H21").Value = 2.062734
$ws.Range("I21").Value = 0.004018499278376935
$ws.Range("J21").Value = 0.004018499278376936
$ws.Range("M21").Value = 63.91341533333334
$ws.Range("N21").Value = 191.740246
$ws.Range("O21").Value = 0.2473779113615065
$ws.Range("P21").Value = 0.2473779113615065
$ws.Range("Q21").Value = 43.94545828806267
$ws.Range("R21").Value = 395.509124592564
$ws.Range("S21").Value = 0.0009940879582926074
$ws.Range("T21").Value = 0.0009940879582926074

# Row 22
$ws.Range("G22").Value = 33.136844
$ws.Range("H22").Value = 99.410532
$ws.Range("I22").Value = 0.1936658585668668
$ws.Range("J22").Value = 0.1936658585668668
$ws.Range("M22").Value = 3.243447333333334
$ws.Range("N22").Value = 9.730342
$ws.Range("O22").Value = 0.01255381554477167
$ws.Range("P22").Value = 0.01255381554477167
$ws.Range("Q22").Value = 107.4776083068827
$ws.Range("R22").Value = 967.298474761944
$ws.Range("S22").Value = 0.002431245465768284
$ws.Range("T22").Value = 0.002431245465768284

# Row 23
$ws.Range("G23").Value = 33.136844
$ws.Range("H23").Value = 99.410532
$ws.Range("I23").Value = 0.1936658585668668
$ws.Range("J23").Value = 0.1936658585668668
$ws.Range("O23").Value = 0.2470858318854151
$ws.Range("P23").Value = 0.2470858318854151
$ws.Range("Q23").Value = 2115.388278794718
$ws.Range("R23").Value = 19038.49450915246
$ws.Range("S23").Value = 0.04785208977179743
$ws.Range("T23").Value = 0.04785208977179743

# Row 24
$ws.Range("G24").Value = 33.136844
$ws.Range("H24").Value = 99.410532
$ws.Range("I24").Value = 0.1936658585668668
$ws.Range("J24").Value = 0.1936658585668668
$ws.Range("M24").Value = 103.7552793333333
$ws.Range("N24").Value = 311.265838
$ws.Range("O24").Value = 0.4015864926064037
$ws.Range("P24").Value = 0.4015864926064037
$ws.Range("Q24").Value = 3438.122505445091
$ws.Range("R24").Value = 30943.10254900581
$ws.Range("S24").Value = 0.07777359287947587
$ws.Range("T24").Value = 0.07777359287947586

# Row 25
$ws.Range("G25").Value = 33.136844
$ws.Range("H25").Value = 99.410532
$ws.Range("I25").Value = 0.1936658585668668
$ws.Range("J25").Value = 0.1936658585668668
$ws.Range("M25").Value = 23.61337433333334
$ws.Range("N25").Value = 70.84012300000001
$ws.Range("O25").Value = 0.09139594860190291
$ws.Range("P25").Value = 0.09139594860190289
$ws.Range("Q25").Value = 782.4727015972709
$ws.Range("R25").Value = 7042.254314375436
$ws.Range("S25").Value = 0.01770027485552075
$ws.Range("T25").Value = 0.01770027485552075

# Row 26
$ws.Range("G26").Value = 33.136844
$ws.Range("H26").Value = 99.410532
$ws.Range("I26").Value = 0.1936658585668668
$ws.Range("J26").Value = 0.1936658585668668
$ws.Range("M26").Value = 63.91341533333334
$ws.Range("N26").Value = 191.740246
$ws.Range("O26").Value = 0.2473779113615065
$ws.Range("P26").Value = 0.2473779113615065
$ws.Range("Q26").Value = 2117.888873407875
$ws.Range("R26").Value = 19060.99986067087
$ws.Range("S26").Value = 0.04790865559430443
$ws.Range("T26").Value = 0.04790865559430442

Write-Output "Updated 279 cell(s) with refreshed TPM-derived values."
